# Refined metadata to be additional tab
#
# 1) Refresh the "time_taken" column (F2:F36) on the existing "data" sheet
#    with the new query timestamps.
# 2) Add a new "metadata" worksheet (after "data") describing the panel
#    query itself (data_name, data_id, data_version, data_version_created,
#    panel_query_time, panel_get_request).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------
# 1) Update the per-gene query timestamps on the "data" sheet.
# ---------------------------------------------------------------------
$newTimestamps = @(
    "2021-10-05 14:35:38.919179",
    "2021-10-05 14:35:38.919187",
    "2021-10-05 14:35:38.919190",
    "2021-10-05 14:35:38.919193",
    "2021-10-05 14:35:38.919196",
    "2021-10-05 14:35:38.919199",
    "2021-10-05 14:35:38.919201",
    "2021-10-05 14:35:38.919204",
    "2021-10-05 14:35:38.919207",
    "2021-10-05 14:35:38.919209",
    "2021-10-05 14:35:38.919212",
    "2021-10-05 14:35:38.919214",
    "2021-10-05 14:35:38.919217",
    "2021-10-05 14:35:38.919219",
    "2021-10-05 14:35:38.919222",
    "2021-10-05 14:35:38.919225",
    "2021-10-05 14:35:38.919228",
    "2021-10-05 14:35:38.919230",
    "2021-10-05 14:35:38.919233",
    "2021-10-05 14:35:38.919235",
    "2021-10-05 14:35:38.919238",
    "2021-10-05 14:35:38.919240",
    "2021-10-05 14:35:38.919243",
    "2021-10-05 14:35:38.919245",
    "2021-10-05 14:35:38.919248",
    "2021-10-05 14:35:38.919251",
    "2021-10-05 14:35:38.919253",
    "2021-10-05 14:35:38.919256",
    "2021-10-05 14:35:38.919258",
    "2021-10-05 14:35:38.919261",
    "2021-10-05 14:35:38.919263",
    "2021-10-05 14:35:38.919266",
    "2021-10-05 14:35:38.919269",
    "2021-10-05 14:35:38.919271",
    "2021-10-05 14:35:38.919274"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# ---------------------------------------------------------------------
# 2) Add the "metadata" sheet right after "data".
# ---------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (bold / centered / bordered, same look as the "data" header).
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Copy the header style from "data"!B1 (bold, centered, thin border) onto
# the new header cells.
$dataSheet.Range("B1").Copy() | Out-Null
$metaSheet.Range("B1:G1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Data row. A2 reuses the bordered/bold "index column" look from data!A2.
$metaSheet.Range("A2").Value = 0
$dataSheet.Range("A2").Copy() | Out-Null
$metaSheet.Range("A2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$metaSheet.Range("B2").Value = "Retinitis pigmentosa_Autosomal Dominant"
$metaSheet.Range("C2").Value = 278

# "data_version" must stay textual ("0.31"), not be coerced into a number -
# format the cell as Text before assigning it, then drop back to the
# default "Normal" style so no stray number-format lingers on the cell.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.31"
$metaSheet.Range("D2").Style = "Normal"

$metaSheet.Range("E2").Value = "2021-08-10T08:07:17.463100Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:38.915535"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/278/?format=json"

$metaSheet.Range("A1").Select() | Out-Null

# Keep "data" as the active/selected sheet, same as before the edit.
$dataSheet.Activate()
$dataSheet.Range("A1").Select() | Out-Null

Write-Output "metadata sheet added; data!F2:F36 timestamps refreshed"
